$d = $word.ActiveDocument

# Locate the "Micro results" row in the (only) table and grab its second cell,
# which holds the long block of micro-result lines that needs to be replaced.
$table = $d.Tables.Item(1)
$targetCell = $null
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $label = $table.Cell($r, 1).Range.Text
    if ($label -like "Micro results*") {
        $targetCell = $table.Cell($r, 2)
        break
    }
}

$endash = [char]0x2013

# Run properties shared by every new paragraph: Times New Roman, blue, 10pt (sz 20 half-points).
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr>'

# Build each new line of text (using [char]0x2013 for the en-dash so the byte
# sequence survives string concatenation correctly).
$line1 = "--------Previous result (1 year)--------"
$line2 = "24/07 " + $endash + " EBV VCA IgG " + $endash + " Positive"
$line3 = "19/06 " + $endash + " RESPIRATORY PCR " + $endash + " **Positive**  "
$line4 = "**Summary:** Human Rhinovirus/Enterovirus detected."
$line5 = "12/06 " + $endash + " EBV VCA IgG " + $endash + " Positive"
$line6 = "20/02 " + $endash + " RESPIRATORY PCR " + $endash + " **Positive**  "
$line7 = "**Summary:** Respiratory Syncytial Virus detected."
$line8 = "04/02 " + $endash + " RESPIRATORY PCR " + $endash + " **Positive**  "
$line9 = "**Summary:** Respiratory Syncytial Virus DETECTED."

# Paragraph 1: a blank paragraph with an empty run (no pPr, no rPr).
$p1 = '<w:p><w:r/></w:p>'

# Paragraph 2: a run with the shared rPr but no text at all.
$p2 = '<w:p><w:r>' + $rPr + '</w:r></w:p>'

# Paragraphs 3-11: one run each carrying text, preserving trailing spaces where present.
$p3  = '<w:p><w:r>' + $rPr + '<w:t>' + $line1 + '</w:t></w:r></w:p>'
$p4  = '<w:p><w:r>' + $rPr + '<w:t>' + $line2 + '</w:t></w:r></w:p>'
$p5  = '<w:p><w:r>' + $rPr + '<w:t xml:space="preserve">' + $line3 + '</w:t></w:r></w:p>'
$p6  = '<w:p><w:r>' + $rPr + '<w:t>' + $line4 + '</w:t></w:r></w:p>'
$p7  = '<w:p><w:r>' + $rPr + '<w:t>' + $line5 + '</w:t></w:r></w:p>'
$p8  = '<w:p><w:r>' + $rPr + '<w:t xml:space="preserve">' + $line6 + '</w:t></w:r></w:p>'
$p9  = '<w:p><w:r>' + $rPr + '<w:t>' + $line7 + '</w:t></w:r></w:p>'
$p10 = '<w:p><w:r>' + $rPr + '<w:t xml:space="preserve">' + $line8 + '</w:t></w:r></w:p>'
$p11 = '<w:p><w:r>' + $rPr + '<w:t>' + $line9 + '</w:t></w:r></w:p>'

$bodyInner = $p1 + $p2 + $p3 + $p4 + $p5 + $p6 + $p7 + $p8 + $p9 + $p10 + $p11

$xmlPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetCell.Range.InsertXML($xmlPayload)
